$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 48; this shifts existing rows 48:59 down to 49:60
$ws.Rows.Item(48).Insert()

# Populate the new row 48 with the new weekly record (same constant columns as
# every other data row, but a new Fecha / Volumen / Precio set).
$ws.Range("A48").Value = 1
$ws.Range("B48").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C48").Value = "Arica y Parinacota"
$ws.Range("D48").Value = 44559
$ws.Range("E48").Value = 15
$ws.Range("F48").Value = 100112038
$ws.Range("G48").Value = "Cebollín baby"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 300
$ws.Range("K48").Value = 3500
$ws.Range("L48").Value = 4000
$ws.Range("M48").Value = 3750
$ws.Range("N48").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O48").Value = "Región de Arica y Parinacota"
$ws.Range("P48").Value = 1875
$ws.Range("Q48").Value = 2
$ws.Range("R48").Value = "Hortaliza"
